# Strip the spurious trailing "16" that was accidentally appended to the
# Bible verse references in column A (e.g. "Zechariah 1:116" -> "Zechariah 1:1"),
# making the reference column human readable again. Column B is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Text
    if ($current -ne $null -and $current -ne "") {
        $fixed = $current -replace "16$", ""
        if ($fixed -ne $current) {
            $cell.Value = $fixed
        }
    }
}
